# CDP Network Audit template update
# - Adds "Serial" and "UPTIME" columns to the Audit sheet header row (between
#   MANAGEMENT_IP and PLATFORM), shifting PLATFORM/SOFTWARE_VERSION/CAPABILITIES
#   two columns to the right.
# - Extends the Audit!_FilterDatabase defined name to cover the two new columns.
# - Populates resolved-hostname / error-IP data collected by the audit run on
#   the "DNS Resolved", "Connection Errors" and "Authentication Errors" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Audit sheet: insert "Serial" / "UPTIME" columns before PLATFORM
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Audit")

$ws1.Columns("G:H").Insert()

# Match the formatting already used across the rest of the header row.
$ws1.Range("F11").Copy()
$ws1.Range("G11:H11").PasteSpecial(-4122)

$ws1.Range("G11").Value = "Serial"
$ws1.Range("H11").Value = "UPTIME"

$ws1.Range("E5").Select() | Out-Null

# ---------------------------------------------------------------------------
# Extend the hidden AutoFilter defined name out to the new last column (K)
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
  if ($n.Name -eq "Audit!_FilterDatabase") {
    $n.RefersTo = "=Audit!`$A`$11:`$K`$11"
  }
}

# ---------------------------------------------------------------------------
# DNS Resolved sheet: reverse-DNS lookup results
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("DNS Resolved")
$ws2.Range("A5").Value = "GB-CAY2-001ASW001"
$ws2.Range("B5").Value = "10.145.61.10"
$ws2.Range("A6").Value = "GB-CAY2-001CSW001"
$ws2.Range("B6").Value = "10.145.63.1"
$ws2.Range("A7").Value = "gb-cay2-001sdw101"
$ws2.Range("B7").Value = "10.255.145.61"
$ws2.Range("A8").Value = "gb-cay2-001sdw102"
$ws2.Range("B8").Value = "10.255.145.62"

# ---------------------------------------------------------------------------
# Connection Errors sheet: IPs that failed to connect
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Connection Errors")
$ws3.Range("A5").Value = "10.250.16.17"
$ws3.Range("A6").Value = "10.250.16.22"
$ws3.Range("A7").Value = "10.250.16.21"
$ws3.Range("A8").Value = "10.250.16.18"

# ---------------------------------------------------------------------------
# Authentication Errors sheet: IPs that failed to authenticate
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Authentication Errors")
$ws4.Range("A5").Value = "62.172.66.138"
